$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same style (bold, bordered, centered)
# as the rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-16
$data = @(
    @(8, 8),
    @(6, 6),
    @(6, 7),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(3, 4),
    @(5, 5),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
